{"js": "// Merge the three adjacent hyperlink runs\n//   `\"\u0423\u0432\u043e\u0434 \u0432 \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u0438\u0440\u0430\u043d\u0435\u0442\u043e\"` + ` ` + `\u0437\u0430 \u0443\u0447\u0435\u043d\u0438\u0446\u0438`\n// into a single run (same displayed text, same \"Hyperlink\" character\n// style) \u2014 i.e. collapse the run split without changing the visible\n// text or the link target.\nconst body = context.document.body;\nconst finalText = '\"\u0423\u0432\u043e\u0434 \u0432 \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u0438\u0440\u0430\u043d\u0435\u0442\u043e\" \u0437\u0430 \u0443\u0447\u0435\u043d\u0438\u0446\u0438';\n\nlet results = body.search(finalText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // The text is already present \u2014 but it may still be split across the\n  // original three runs. `insertText('Replace')` is a no-op when the\n  // replacement text is byte-identical to what's already there, so it\n  // would not actually coalesce the runs. Route the replace through a\n  // throw-away placeholder first: that guarantees a real text change,\n  // which forces Word to collapse the target span into one run; the\n  // second replace then restores the original text in that single run.\n  const placeholder = \"\\u0001__TMP_MERGE_PLACEHOLDER__\\u0001\";\n  results.items[0].insertText(placeholder, \"Replace\");\n  await context.sync();\n\n  const ph = body.search(placeholder, { matchCase: true });\n  ph.load(\"items\");\n  await context.sync();\n\n  const merged = ph.items[0].insertText(finalText, \"Replace\");\n  // Re-apply the hyperlink's character style (\"aa\" -> built-in\n  // \"Hyperlink\") to the newly-created run, since a plain text replace\n  // drops direct/style formatting that used to live on the replaced runs.\n  merged.style = \"Hyperlink\";\n  await context.sync();\n}\n", "ps1": "# Merge the three adjacent hyperlink runs\n#   `\"\u0423\u0432\u043e\u0434 \u0432 \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u0438\u0440\u0430\u043d\u0435\u0442\u043e\"` + ` ` + `\u0437\u0430 \u0443\u0447\u0435\u043d\u0438\u0446\u0438`\n# into a single run (same displayed text, same \"Hyperlink\" character\n# style) \u2014 i.e. collapse the run split without changing the visible\n# text or the link target.\n\n$d = $word.ActiveDocument\n$finalText = '\"\u0423\u0432\u043e\u0434 \u0432 \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u0438\u0440\u0430\u043d\u0435\u0442\u043e\" \u0437\u0430 \u0443\u0447\u0435\u043d\u0438\u0446\u0438'\n\n# Locate the hyperlink by its (distinctive) display text rather than\n# assuming it is the first one in the document.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {\n  $candidate = $d.Hyperlinks.Item($i)\n  if ($candidate.Range.Text -like '*\u0437\u0430 \u0443\u0447\u0435\u043d\u0438\u0446\u0438*') {\n    $targetIndex = $i\n  }\n}\n\nif ($targetIndex -ge 1) {\n  $hyperlink = $d.Hyperlinks.Item($targetIndex)\n  $start = $hyperlink.Range.Start\n  $end = $hyperlink.Range.End\n\n  # Re-seat onto a plain Document.Range (using the raw Hyperlink.Range\n  # object for the write below is unreliable), then force a genuine text\n  # change first: setting identical text back is a no-op and would leave\n  # the original run split intact, so route the replace through a\n  # throw-away placeholder. That guarantees Word collapses the run span\n  # into a single run; the second assignment restores the exact original\n  # text in that single run.\n  $range = $d.Range($start, $end)\n  $range.Text = \"TmpMergePlaceholder\"\n\n  $hyperlink2 = $d.Hyperlinks.Item($targetIndex)\n  $start2 = $hyperlink2.Range.Start\n  $end2 = $hyperlink2.Range.End\n  $range2 = $d.Range($start2, $end2)\n  $range2.Text = $finalText\n\n  # Re-apply the hyperlink's character style (\"aa\" -> built-in\n  # \"Hyperlink\") to the newly-created run, since the plain text\n  # assignment drops direct/style formatting that used to live on the\n  # replaced runs.\n  $hyperlink3 = $d.Hyperlinks.Item($targetIndex)\n  $start3 = $hyperlink3.Range.Start\n  $end3 = $hyperlink3.Range.End\n  $range3 = $d.Range($start3, $end3)\n  $range3.Style = \"Hyperlink\"\n}\n"}
